$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "disable" column (D) with header + Yes/No values
$ws.Range("D1").Value = "disable"
$ws.Range("D2").Value = "Yes"
$ws.Range("D3").Value = "No"
$ws.Range("D4").Value = "Yes"

# B4 used to hold text "NC"; it becomes a real date (10/22/2014) formatted MM/DD/YY
$ws.Range("B4").Value = 41934
$ws.Range("B4").NumberFormat = "MM/DD/YY"

# Match the author's final selection
$ws.Range("B5").Select()
